# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (fund-holdings detail, like the other
#    quarterly sheets) right after "2021-Q4" and before "总计".
# 2. Insert a new summary row for "2022-Q1" at the top of the "总计" sheet's
#    data block (pushing the existing rows down by one and renumbering the
#    leading index column).

$wb = $excel.ActiveWorkbook

# Helper: write a row of values starting at column B (col 2) through a
# scratch range far to the right, then paste-special VALUES ONLY into the
# destination. Routing through a scratch cell/range that was populated with
# a leading single-quote forces Excel to store numeric-looking strings
# ("007139", "12.79", ...) as text instead of auto-converting them to
# numbers; pasting only the *values* afterwards drops the transient
# quote-prefix formatting it picks up along the way, so the destination
# cells end up with no stray style applied.
function Set-RowValues {
    param($sheet, $rowNum, [object[]]$values)

    $n = $values.Count
    $arr = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) {
        $v = $values[$i]
        if ($v -is [string]) {
            $arr[0, $i] = "'" + $v
        } else {
            $arr[0, $i] = $v
        }
    }

    $scratch = $sheet.Range($sheet.Cells.Item($rowNum, 26), $sheet.Cells.Item($rowNum, 25 + $n))
    $scratch.Value = $arr
    $scratch.Copy()
    $dest = $sheet.Range($sheet.Cells.Item($rowNum, 2), $sheet.Cells.Item($rowNum, 1 + $n))
    $dest.PasteSpecial(-4163)
    $scratch.Clear()
}

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet
# ---------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $srcSheet)
$newSheet.Name = "2022-Q1"

# Reuse the existing header / index-column formatting (bold, centered,
# bordered) from the "2021-Q4" sheet so the new sheet matches the others.
$srcSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$srcSheet.Range("A2:A7").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)
$srcSheet.Range("A7").Copy()
$newSheet.Range("A8").PasteSpecial(-4122)

for ($i = 0; $i -le 6; $i++) {
    $newSheet.Cells.Item(2 + $i, 1).Value = $i
}

Set-RowValues $newSheet 1 @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
Set-RowValues $newSheet 2 @("007139", "富国民裕进取沪港深成长精选混合", "12.79", "92.21", "5.60", "0.7162", 5)
Set-RowValues $newSheet 3 @("005504", "汇添富沪港深大盘价值混合", "3.49", "92.33", "4.81", "0.1679", 9)
Set-RowValues $newSheet 4 @("015119", "汇添富沪港深大盘价值混合D", "3.49", "92.33", "4.81", "0.1679", 9)
Set-RowValues $newSheet 5 @("006205", "汇添富沪港深优势精选定期开放混合", "0.40", "93.67", "6.53", "0.0261", 3)
Set-RowValues $newSheet 6 @("005142", "中融沪港深大消费主题灵活配置混合A", "0.49", "88.98", "5.31", "0.0260", 5)
Set-RowValues $newSheet 7 @("005143", "中融沪港深大消费主题灵活配置混合C", "0.33", "88.98", "5.31", "0.0175", 5)
Set-RowValues $newSheet 8 @("005269", "华泰柏瑞港股通量化灵活配置混合", "0.33", "37.77", "0.87", "0.0029", 7)

# ---------------------------------------------------------------------
# 2) "总计" sheet: insert the 2022-Q1 summary row at the top of the table
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

# Re-apply the bordered/centered index-column style (copied from the row
# below, which still carries it) to the freshly inserted A2 cell.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 1.12

for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
